# "Generate Report for Handback" - mark the two handed-off files as handed
# back (in sync with en-US), stamp new handback datetime values for the
# zh-cn and de-de targets, and point the "Latest Target File"/"Latest
# Handback File" columns at the appropriate files (with a hyperlink on the
# target-file cell), widening a couple of columns to fit the new content.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: status column (E/F) for both rows now reads "Handed
# back: in sync with en-US" instead of "Ready for handoff". Because both
# cells in both rows share the same underlying text, setting all four
# keeps the workbook's string table de-duplicated exactly like Excel does.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# Widen the two status columns so the longer text fits.
$wsOverview.Columns("E").ColumnWidth = 29.166666666666668
$wsOverview.Columns("F").ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# "Status" column also shows the same handed-back text.
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

# Latest Target File (I) now references the source .md file and is a
# hyperlink to it, same as column A.
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6bbe59e3137bda7925152250914f5dd084fd197b/e2e/110ee19f-930f-4e70-b45c-16ab885fc797.md", [Type]::Missing, [Type]::Missing, "110ee19f-930f-4e70-b45c-16ab885fc797.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6bbe59e3137bda7925152250914f5dd084fd197b/e2e/c086225b-de91-444b-956a-b10c7451926b.md", [Type]::Missing, [Type]::Missing, "c086225b-de91-444b-956a-b10c7451926b.md") | Out-Null

# Latest Handback File (J) now points at the generated xlf handoff file
# for each row (matching column G, the Latest Handoff File).
$wsZh.Range("J2").Value = $wsZh.Range("G2").Text
$wsZh.Range("J3").Value = $wsZh.Range("G3").Text

# Latest Handback DateTime (K) gets stamped with the handback run time.
$wsZh.Range("K2").Value = "2016-09-06 08:39:02"
$wsZh.Range("K3").Value = "2016-09-06 08:39:02"

# Widen the columns that now hold longer filenames.
$wsZh.Columns("C").ColumnWidth = 29.166666666666668
$wsZh.Columns("I").ColumnWidth = 39.166666666666664
$wsZh.Columns("J").ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet (same shape as zh-cn, different target language files)
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6bbe59e3137bda7925152250914f5dd084fd197b/e2e/110ee19f-930f-4e70-b45c-16ab885fc797.md", [Type]::Missing, [Type]::Missing, "110ee19f-930f-4e70-b45c-16ab885fc797.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6bbe59e3137bda7925152250914f5dd084fd197b/e2e/c086225b-de91-444b-956a-b10c7451926b.md", [Type]::Missing, [Type]::Missing, "c086225b-de91-444b-956a-b10c7451926b.md") | Out-Null

$wsDe.Range("J2").Value = $wsDe.Range("G2").Text
$wsDe.Range("J3").Value = $wsDe.Range("G3").Text

# de-de's handback run happened a bit later than zh-cn's.
$wsDe.Range("K2").Value = "2016-09-06 08:39:21"
$wsDe.Range("K3").Value = "2016-09-06 08:39:21"

$wsDe.Columns("C").ColumnWidth = 29.166666666666668
$wsDe.Columns("I").ColumnWidth = 39.166666666666664
$wsDe.Columns("J").ColumnWidth = 39.166666666666664
